$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving numeric-looking values must be forced to Text format
# so Excel stores/reads them back as strings (matching the source data),
# instead of auto-converting them to numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated cryptos list values
$ws.Range("D2").Value = '96.570.43'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '3.729.00'
$ws.Range("E3").Value = '  +3.72%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '238.83'
$ws.Range("E5").Value = '  -2.26%  '
$ws.Range("D6").Value = '1.93'
$ws.Range("E6").Value = '  +8.33%  '
$ws.Range("D7").Value = '656.53'
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("D8").Value = '0.423'
$ws.Range("E8").Value = '  -0.76%  '
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("D10").Value = '0.999'
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("D11").Value = '3.733.09'
$ws.Range("E11").Value = '  +3.92%  '
$ws.Range("D12").Value = '45.19'
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("E13").Value = '  +0.57%  '
$ws.Range("D14").Value = '6.87'
$ws.Range("E14").Value = '  +6.18%  '
$ws.Range("D15").Value = '4.421.16'
$ws.Range("E15").Value = '  +3.55%  '
$ws.Range("D16").Value = '0.0000269'
$ws.Range("E16").Value = '  +3.13%  '
$ws.Range("D17").Value = '96.403.90'
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").Value = '9.00'
$ws.Range("E18").Value = '  +16.13%  '
$ws.Range("D19").Value = '3.732.81'
$ws.Range("E19").Value = '  +3.67%  '
$ws.Range("D20").Value = '19.11'
$ws.Range("E20").Value = '  +4.26%  '
$ws.Range("D21").Value = '12.85'
$ws.Range("E21").Value = '  +1.26%  '
$ws.Range("D22").Value = '0.527'
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("D23").Value = '523.86'
$ws.Range("E23").Value = '  +1.09%  '
$ws.Range("D24").Value = '3.49'
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("D25").Value = '7.07'
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("D26").Value = '0.0000204'
$ws.Range("E26").Value = '  -1.57%  '
$ws.Range("D27").Value = '102.02'
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("D28").Value = '13.41'
$ws.Range("E28").Value = '  +0.73%  '
$ws.Range("E29").Value = '  -8.45%  '
$ws.Range("D30").Value = '12.51'
$ws.Range("E30").Value = '  +3.98%  '
$ws.Range("D31").Value = '3.07'
$ws.Range("E31").Value = '  +3.02%  '
$ws.Range("E32").Value = '  +0.19%  '
$ws.Range("E33").Value = '  +12.09%  '
$ws.Range("D34").Value = '0.186'
$ws.Range("E34").Value = '  -2.37%  '
$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").Value = '669.83'
$ws.Range("E35").Value = '  +8.92%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '32.92'
$ws.Range("E36").Value = '  +3.29%  '
$ws.Range("B37").Value = 'Binance-PegBSC-USD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("D38").Value = '0.598'
$ws.Range("E38").Value = '  +2.09%  '
$ws.Range("D39").Value = '8.90'
$ws.Range("E39").Value = '  +1.68%  '
$ws.Range("D40").Value = '7.07'
$ws.Range("E40").Value = '  +15.95%  '
$ws.Range("D41").Value = '41.15'
$ws.Range("E41").Value = '  +24.58%  '
$ws.Range("E42").Value = '  +4.07%  '
$ws.Range("D43").Value = '0.984'
$ws.Range("E43").Value = '  +5.79%  '
$ws.Range("D44").Value = '1.99'
$ws.Range("E44").Value = '  +3.07%  '
$ws.Range("D46").Value = '0.446'
$ws.Range("E46").Value = '  -3.61%  '
$ws.Range("D47").Value = '0.0458'
$ws.Range("E47").Value = '  +2.35%  '
$ws.Range("D48").Value = '2.34'
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("D49").Value = '23.62'
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("D50").Value = '8.60'
$ws.Range("E50").Value = '  -0.94%  '
$ws.Range("E51").Value = '  +2.41%  '
